$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.883886666666666
$ws.Range("H2").Value = 5.65166
$ws.Range("I2").Value = 0.7298568945019562
$ws.Range("J2").Value = 0.7298568945019563
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 17.16653
$ws.Range("N2").Value = 51.49959
$ws.Range("O2").Value = 0.0560345397128279
$ws.Range("P2").Value = 0.0560345397128279
$ws.Range("Q2").Value = 32.33979697993333
$ws.Range("R2").Value = 291.0581728194
$ws.Range("S2").Value = 0.04089719513965111
$ws.Range("T2").Value = 0.04089719513965111
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.883886666666666
$ws.Range("H3").Value = 5.65166
$ws.Range("I3").Value = 0.7298568945019562
$ws.Range("J3").Value = 0.7298568945019563
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 256.4443053333333
$ws.Range("N3").Value = 769.332916
$ws.Range("O3").Value = 0.8370788162388805
$ws.Range("P3").Value = 0.8370788162388805
$ws.Range("Q3").Value = 483.1120075600621
$ws.Range("R3").Value = 4348.00806804056
$ws.Range("S3").Value = 0.610947745273483
$ws.Range("T3").Value = 0.6109477452734831
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.883886666666666
$ws.Range("H4").Value = 5.65166
$ws.Range("I4").Value = 0.7298568945019562
$ws.Range("J4").Value = 0.7298568945019563
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 32.74538866666666
$ws.Range("N4").Value = 98.236166
$ws.Range("O4").Value = 0.1068866440482915
$ws.Range("P4").Value = 0.1068866440482915
$ws.Range("Q4").Value = 61.6886011039511
$ws.Range("R4").Value = 555.1974099355599
$ws.Range("S4").Value = 0.07801195408882206
$ws.Range("T4").Value = 0.07801195408882207
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.353459
$ws.Range("H5").Value = 1.060377
$ws.Range("I5").Value = 0.1369373713601492
$ws.Range("J5").Value = 0.1369373713601492
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 17.16653
$ws.Range("N5").Value = 51.49959
$ws.Range("O5").Value = 0.0560345397128279
$ws.Range("P5").Value = 0.0560345397128279
$ws.Range("Q5").Value = 6.06766452727
$ws.Range("R5").Value = 54.60898074543
$ws.Range("S5").Value = 0.007673222573650542
$ws.Range("T5").Value = 0.007673222573650544
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.353459
$ws.Range("H6").Value = 1.060377
$ws.Range("I6").Value = 0.1369373713601492
$ws.Range("J6").Value = 0.1369373713601492
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 256.4443053333333
$ws.Range("N6").Value = 769.332916
$ws.Range("O6").Value = 0.8370788162388805
$ws.Range("P6").Value = 0.8370788162388805
$ws.Range("Q6").Value = 90.64254771881467
$ws.Range("R6").Value = 815.7829294693321
$ws.Range("S6").Value = 0.1146273727170177
$ws.Range("T6").Value = 0.1146273727170177
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.353459
$ws.Range("H7").Value = 1.060377
$ws.Range("I7").Value = 0.1369373713601492
$ws.Range("J7").Value = 0.1369373713601492
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 32.74538866666666
$ws.Range("N7").Value = 98.236166
$ws.Range("O7").Value = 0.1068866440482915
$ws.Range("P7").Value = 0.1068866440482915
$ws.Range("Q7").Value = 11.57415233273133
$ws.Range("R7").Value = 104.167370994582
$ws.Range("S7").Value = 0.01463677606948098
$ws.Range("T7").Value = 0.01463677606948098
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.343827
$ws.Range("H8").Value = 1.031481
$ws.Range("I8").Value = 0.1332057341378944
$ws.Range("J8").Value = 0.1332057341378944
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 17.16653
$ws.Range("N8").Value = 51.49959
$ws.Range("O8").Value = 0.0560345397128279
$ws.Range("P8").Value = 0.0560345397128279
$ws.Range("Q8").Value = 5.90231651031
$ws.Range("R8").Value = 53.12084859279
$ws.Range("S8").Value = 0.00746412199952624
$ws.Range("T8").Value = 0.00746412199952624
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.343827
$ws.Range("H9").Value = 1.031481
$ws.Range("I9").Value = 0.1332057341378944
$ws.Range("J9").Value = 0.1332057341378944
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 256.4443053333333
$ws.Range("N9").Value = 769.332916
$ws.Range("O9").Value = 0.8370788162388805
$ws.Range("P9").Value = 0.8370788162388805
$ws.Range("Q9").Value = 88.172476169844
$ws.Range("R9").Value = 793.552285528596
$ws.Range("S9").Value = 0.1115036982483797
$ws.Range("T9").Value = 0.1115036982483797
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.343827
$ws.Range("H10").Value = 1.031481
$ws.Range("I10").Value = 0.1332057341378944
$ws.Range("J10").Value = 0.1332057341378944
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 32.74538866666666
$ws.Range("N10").Value = 98.236166
$ws.Range("O10").Value = 0.1068866440482915
$ws.Range("P10").Value = 0.1068866440482915
$ws.Range("Q10").Value = 11.258748749094
$ws.Range("R10").Value = 101.328738741846
$ws.Range("S10").Value = 0.01423791388998848
$ws.Range("T10").Value = 0.01423791388998848